$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

for ($r = 2; $r -le 319; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
